$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Straightforward cell value updates (text stays text; Excel auto-numeric-parses
# only where that matches the source diff's intent).
$ws.Range("D2").Value = "29.422.03"
$ws.Range("D3").Value = "1.872.84"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "0.7114"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("D6").Value = "241.98"
$ws.Range("E6").Value = "  +1.69%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.07905"
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("D9").Value = "0.3112"
$ws.Range("E9").Value = "  +2.89%  "
$ws.Range("E10").Value = "  +5.94%  "
$ws.Range("D11").Value = "0.08259"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D12").Value = "1.882.77"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.7256"
$ws.Range("E13").Value = "  +2.64%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.281"
$ws.Range("E14").Value = "  +1.50%  "
$ws.Range("D15").Value = "90.84"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("D16").Value = "29.447.86"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("E17").Value = "  +1.69%  "
$ws.Range("D18").Value = "247.41"
$ws.Range("E18").Value = "  +4.54%  "
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.131.07"
$ws.Range("E21").Value = "  +2.03%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("E23").Value = "  +6.21%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").Value = "0.1585"
$ws.Range("E25").Value = "  +12.48%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "163.67"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "9.006"
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("E29").Value = "  -3.41%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "1.496"
$ws.Range("E30").Value = "  +1.73%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "4.377"
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "4.117"
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.05306"
$ws.Range("E33").Value = "  +2.06%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.195"
$ws.Range("E35").Value = "  +2.64%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.7229"
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.679"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01868"
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.238.86"
$ws.Range("E39").Value = "  +7.36%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.719"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "0.9104"
$ws.Range("E41").Value = "  -2.61%  "
$ws.Range("D42").Value = "6.188"
$ws.Range("E42").Value = "  +3.27%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "73.92"
$ws.Range("E43").Value = "  +5.17%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "102.88"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "0.5328"
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "2.029.63"
$ws.Range("E47").Value = "  +2.53%  "
$ws.Range("B48").Value = "SynthetixNetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D48").Value = "2.947"
$ws.Range("E48").Value = "  +13.47%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.759"
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "9.293"
$ws.Range("E50").Value = "  +1.58%  "
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").Value = "0.4314"
$ws.Range("E51").Value = "  +1.18%  "

# These price cells must stay literal text (e.g. trailing/insignificant zeros
# like "1.000" or "18.30") -- a plain .Value assignment would let Excel
# reinterpret them as numbers and silently normalize the digits away.
# Flip the cell to Text format, write the string, then drop the format override
# again so the cell keeps its original (unstyled) look.
$textCells = @("D19", "D23", "D28", "D29", "D34", "D44")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}
$ws.Range("D19").Value = "0.000007864"
$ws.Range("D23").Value = "7.940"
$ws.Range("D28").Value = "18.30"
$ws.Range("D29").Value = "1.360"
$ws.Range("D34").Value = "1.930"
$ws.Range("D44").Value = "1.000"
foreach ($cell in $textCells) {
    $ws.Range($cell).ClearFormats()
}

